$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark (it sat at the end of the
#    "Module 2 - Precision Turning (Advanced)" paragraph before this
#    edit; Word re-creates it at the new last-edit location on save).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Find the "Module 3 - " paragraph.
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "Module 3*") {
        $target = $para
        break
    }
}

$r = $target.Range

# Collapse to just before the paragraph mark (end of the visible text).
$r.MoveEnd(1, -1)
$r.Collapse(0)
$paraStart = $target.Range.Start

$addedText = "Precision Forward Movement"
$prefix = "Precision Forward "
$prefixLen = $prefix.Length

# Typing the new text through the Range.Text setter merges it with the
# existing formatted run, so it correctly inherits the surrounding
# run's formatting (<w:lang w:val="en-US"/>) instead of coming out as
# a bare, unformatted run.
$r.Text = $addedText

# The paragraph is now a single run: "Module 3 - Precision Forward Movement"
$b1 = $paraStart + 8                 # right after "Module 3"
$b2 = $b1 + 3                        # right after " - "
$b3 = $b2 + $prefixLen               # right after "Precision Forward "

# Re-establish the original run boundaries ("Module 3" | " - ") and
# split out the new "Precision Forward " / "Movement" runs by doing a
# self-assignment through FormattedText - this forces the engine to
# materialize a distinct <w:r> at each boundary while preserving the
# (identical) run formatting on both sides of the split.
$seg1 = $d.Range($b2, $b3)
$seg1.FormattedText = $seg1.Duplicate.FormattedText

$seg2 = $d.Range($paraStart, $b1)
$seg2.FormattedText = $seg2.Duplicate.FormattedText

# ------------------------------------------------------------------
# 3. Re-insert the "_GoBack" bookmark between "Precision Forward " and
#    "Movement", matching the new last-edit location.
# ------------------------------------------------------------------
$bmRange = $d.Range($b3, $b3)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "OK"
